$d = $word.ActiveDocument
$sec = $d.Sections(1)
$hdr = $sec.Headers(2)

# delete the inline shape's run
$shp = $hdr.Range.InlineShapes(1)
$shp.Delete()

# now delete whatever remains in the header paragraph (the 'rtl' empty run)
$rest = $hdr.Range
Write-Host ("rest start=" + $rest.Start + " end=" + $rest.End + " text=[" + $rest.Text + "]")
$rest.Delete()
Write-Host ("after rest delete: header xml:")
Write-Host $hdr.Range.WordOpenXML
